$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 275.9
$ws.Range("I11").Value = 275.9
$ws.Range("K11").Value = 275.9
$ws.Range("M11").Value = -135.9
$ws.Range("H137").Value = 1458.6
$ws.Range("J137").Value = 2263
$ws.Range("L137").Value = 6789
$ws.Range("N137").Value = -11889

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 300
$ws.Range("I4").Value = 300
$ws.Range("K4").Value = 300
$ws.Range("M4").Value = -184
$ws.Range("H28").Value = 9957
$ws.Range("I28").Value = 9957
$ws.Range("K28").Value = 9957
$ws.Range("M28").Value = -9765
$ws.Range("H45").Value = 2144
$ws.Range("I45").Value = 2124.75
$ws.Range("K45").Value = 2124.75
$ws.Range("M45").Value = -1747.75
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H99").Value = 9957
$ws.Range("I99").Value = 9957
$ws.Range("K99").Value = 9957
$ws.Range("M99").Value = -6962
$ws.Range("H122").Value = 1019.6875
$ws.Range("I122").Value = 1019.6875
$ws.Range("K122").Value = 3059.0625
$ws.Range("M122").Value = -609.0625
$ws.Range("H132").Value = 1879.1562
$ws.Range("I132").Value = 1573.9231
$ws.Range("K132").Value = 4721.7693
$ws.Range("M132").Value = -2191.7693

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 30
$ws.Range("I11").Value = 30
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 30
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 110
$ws.Range("N11").ClearContents()
$ws.Range("H20").Value = 5544.154
$ws.Range("I20").Value = 3619.4443
$ws.Range("J20").Value = 9874.75
$ws.Range("K20").Value = 3619.4443
$ws.Range("L20").Value = 9874.75
$ws.Range("M20").Value = -3372.4443
$ws.Range("N20").Value = -10368.75
$ws.Range("H99").Value = 2325
$ws.Range("I99").Value = 1650
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 1650
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -152
$ws.Range("N99").Value = -5996

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3955.8333
$ws.Range("J31").Value = 4999.5
$ws.Range("L31").Value = 4999.5
$ws.Range("N31").Value = -5589.5
$ws.Range("H34").Value = 3955.8333
$ws.Range("J34").Value = 4999.5
$ws.Range("L34").Value = 4999.5
$ws.Range("N34").Value = -5403.5
$ws.Range("H58").Value = 11637.333
$ws.Range("J58").Value = 13342.333
$ws.Range("L58").Value = 13342.333
$ws.Range("N58").Value = -13748.333
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H99").Value = 7004.4375
$ws.Range("I99").Value = 6395
$ws.Range("J99").Value = 8345.200000000001
$ws.Range("K99").Value = 6395
$ws.Range("L99").Value = 8345.200000000001
$ws.Range("M99").Value = -4897
$ws.Range("N99").Value = -11341.2
$ws.Range("H126").Value = 7004.4375
$ws.Range("I126").Value = 6395
$ws.Range("J126").Value = 8345.200000000001
$ws.Range("K126").Value = 19185
$ws.Range("L126").Value = 25035.6
$ws.Range("M126").Value = -16715
$ws.Range("N126").Value = -29975.6
$ws.Range("H134").Value = 2403.8
$ws.Range("I134").Value = 2004.4445
$ws.Range("J134").Value = 5998
$ws.Range("K134").Value = 6013.333500000001
$ws.Range("L134").Value = 17994
$ws.Range("M134").Value = -3478.333500000001
$ws.Range("N134").Value = -23064
$ws.Range("H136").Value = 11637.333
$ws.Range("J136").Value = 13342.333
$ws.Range("L136").Value = 40026.999
$ws.Range("N136").Value = -45126.999
$ws.Range("H138").Value = 5561.625
$ws.Range("J138").Value = 34000
$ws.Range("L138").Value = 34000
$ws.Range("N138").Value = -44280
$ws.Range("H141").Value = 356431.7
$ws.Range("J141").Value = 390479.66
$ws.Range("L141").Value = 390479.66
$ws.Range("N141").Value = -400839.66

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 95.210526
$ws.Range("I2").Value = 95.14286
$ws.Range("J2").Value = 95.25
$ws.Range("K2").Value = 570.85716
$ws.Range("L2").Value = 571.5
$ws.Range("M2").Value = -457.85716
$ws.Range("N2").Value = -797.5
$ws.Range("H134").Value = 55560550
$ws.Range("I134").Value = 55560550
$ws.Range("K134").Value = 166681650
$ws.Range("M134").Value = -166676580

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 17996.111
$ws.Range("J15").Value = 17996.111
$ws.Range("L15").Value = 17996.111
$ws.Range("N15").Value = -18572.111
$ws.Range("H70").Value = 10008974
$ws.Range("I70").Value = 14296242
$ws.Range("J70").Value = 5348.6665
$ws.Range("K70").Value = 14296242
$ws.Range("L70").Value = 5348.6665
$ws.Range("M70").Value = -14295972
$ws.Range("N70").Value = -5888.6665
$ws.Range("H73").Value = 10008974
$ws.Range("I73").Value = 14296242
$ws.Range("J73").Value = 5348.6665
$ws.Range("K73").Value = 14296242
$ws.Range("L73").Value = 5348.6665
$ws.Range("M73").Value = -14295306
$ws.Range("N73").Value = -7220.6665
$ws.Range("H81").Value = 17996.111
$ws.Range("J81").Value = 17996.111
$ws.Range("L81").Value = 17996.111
$ws.Range("N81").Value = -19992.111
$ws.Range("H84").Value = 17996.111
$ws.Range("J84").Value = 17996.111
$ws.Range("L84").Value = 53988.333
$ws.Range("N84").Value = -63972.333
$ws.Range("H132").Value = 3001.6316
$ws.Range("I132").Value = 2478.3333
$ws.Range("J132").Value = 3898.7144
$ws.Range("K132").Value = 7434.999899999999
$ws.Range("L132").Value = 11696.1432
$ws.Range("M132").Value = -4904.999899999999
$ws.Range("N132").Value = -16756.1432

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H99").Value = 90259
$ws.Range("I99").Value = 90259
$ws.Range("K99").Value = 90259
$ws.Range("M99").Value = -87264
$ws.Range("H122").Value = 2682.5
$ws.Range("I122").Value = 2682.5
$ws.Range("K122").Value = 8047.5
$ws.Range("M122").Value = -5597.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").ClearContents()
$ws.Range("H122").Value = 1597
$ws.Range("I122").Value = 1597
$ws.Range("K122").Value = 4791
$ws.Range("M122").Value = -2341
$ws.Range("H132").Value = 6500
$ws.Range("J132").Value = 14000
$ws.Range("L132").Value = 42000
$ws.Range("N132").Value = -47060
